$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (names, URLs, percentage strings, multi-dot numbers -
# Excel keeps these as text automatically).
$ws.Range("D2").Value = "68.910.84"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "2.665.02"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("E6").Value = "  +4.98%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "2.661.45"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("E10").Value = "  +14.34%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("E15").Value = "  +6.77%  "
$ws.Range("D16").Value = "3.147.03"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "68.774.04"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "2.663.71"
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  +5.91%  "
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("E28").Value = "  +9.72%  "
$ws.Range("D29").Value = "2.808.78"
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  +5.35%  "
$ws.Range("E33").Value = "  +6.09%  "
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  +5.71%  "
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("E41").Value = "  +6.13%  "
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("E43").Value = "  +8.44%  "
$ws.Range("E44").Value = "  +5.97%  "
$ws.Range("D45").Value = "0.0₆0323"
$ws.Range("E45").Value = "  +15.19%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("E51").Value = "  +4.77%  "

# Numeric-looking price text must be forced to stay text (matches the
# original inline-string cells) - otherwise Excel auto-converts them to
# numbers and mangles trailing zeros / precision.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.88"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.02"
$ws.Range("D6").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.13"
$ws.Range("D14").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.45"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.75"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.46"
$ws.Range("D21").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.90"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.95"
$ws.Range("D25").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("D27").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "580.54"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.42"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.01"
$ws.Range("D33").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.43"
$ws.Range("D38").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.33"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.44"
$ws.Range("D41").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.80"
$ws.Range("D48").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.12"
$ws.Range("D51").ClearFormats()
